# Update on 20181115.2017by YKBKyle
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data (File Management / read.csv())
# Shared-string insertion order must match: read.csv(), read.csv("csv file",as.is=TRUE),
# description, then File Management -- so populate B/C/D before A.
$ws.Range("B36").Value = "read.csv()"
$ws.Range("C36").Value = 'read.csv("csv file",as.is=TRUE)'
$ws.Range("D36").Value = "read csv file, and prevent character from being converted to factor class"
$ws.Range("A36").Value = "File Management"

# Scroll / selection state to match author's final view
$ws.Application.ActiveWindow.ScrollRow = 19
$ws.Range("A37").Select()
